# Confirm the submission date: insert "Thu 26th " (with superscript "th")
# right after "Due Date: " and move the _GoBack bookmark to mark this as
# the most recent edit location (matching Word's own behaviour of
# relocating the single _GoBack bookmark to the latest edit point).

$d = $word.ActiveDocument

# Locate "Due Date: " (unique in the document) and collapse to its end.
$rng = $d.Content
$found = $rng.Find.Execute("Due Date: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Due Date: ' in the document"
}
$rng.Collapse(0)
$insertStart = $rng.Start

# Insert the whole new phrase as plain text first.
$rng.InsertAfter("Thu 26th ")

# Make just the "th" (2 chars, right after "Thu 26") superscript, to match
# "Thu 26[th superscript] ".
$supRange = $d.Range($insertStart + 6, $insertStart + 8)
$supRange.Font.Superscript = $true

# Drop a zero-length range right after the trailing space and (re)plant the
# _GoBack bookmark there -- Word keeps only one _GoBack bookmark per
# document and automatically relocates it to the newest edit, removing it
# from its previous location.
$bmRange = $d.Range($insertStart + 9, $insertStart + 9)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "Inserted date confirmation and relocated _GoBack bookmark."
